$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Summary sheet: update capital / P&L / trade-count stats after the
# new trade (#9) closed.
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.88   # Current Capital
$summary.Range("B4").Value = -0.12     # Total P&L $
$summary.Range("B5").Value = -0.27     # Total P&L %
$summary.Range("B6").Value = 9         # Total Trades
$summary.Range("B8").Value = 4         # Losing Trades
$summary.Range("B9").Value = 33.33     # Win Rate %

# ------------------------------------------------------------------
# Strategy Status sheet: update the MarketMaking strategy row (row 4)
# to reflect the new trade.
# ------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.88      # Capital
$status.Range("D4").Value = 9          # Trades
$status.Range("E4").Value = -0.12      # P&L $
$status.Range("F4").Value = -0.12      # P&L %
$status.Range("G4").Value = 33.33      # Win Rate %

# ------------------------------------------------------------------
# Append the newly closed trade (#9) to both the "All Trades" and the
# "MarketMaking" trade logs as row 10.
# ------------------------------------------------------------------
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A10").Value = 9

    # Force the date to be stored as literal text instead of letting
    # Excel auto-convert "2026-02-17" into a date serial number, then
    # reset the style back to Normal so no stray per-cell formatting
    # is left behind.
    $ws.Range("B10").NumberFormat = "@"
    $ws.Range("B10").Value = "2026-02-17"
    $ws.Range("B10").Style = "Normal"

    $ws.Range("C10").Value = "15:14:13"
    $ws.Range("D10").Value = "MarketMaking"
    $ws.Range("E10").Value = "DOWN"
    $ws.Range("F10").Value = 0.05
    $ws.Range("G10").Value = 0.03
    $ws.Range("H10").Value = "CLOSED"
    $ws.Range("I10").Value = -40
    $ws.Range("J10").Value = -0.02
    $ws.Range("K10").Value = 99.88
    $ws.Range("L10").Value = 0
    $ws.Range("M10").Value = 0
    $ws.Range("N10").Value = 0.6
    $ws.Range("O10").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P10").Value = "early_exit"
    $ws.Range("Q10").Value = 0.14
}
